# Apply the data corrections described in the diff.
# The workbook contains a single table of enrollment/registration numbers;
# several "Inscritos" (E), "Pagos" (F) and "Inscrições homologadas" (H)
# counts were incremented by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 478
$ws.Range("E20").Value = 83
$ws.Range("E27").Value = 311
$ws.Range("E29").Value = 162
$ws.Range("E31").Value = 70

$ws.Range("E33").Value = 274
$ws.Range("F33").Value = 140
$ws.Range("H33").Value = 140

$ws.Range("E35").Value = 143
$ws.Range("E42").Value = 354
$ws.Range("E47").Value = 428
